$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.584.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.484.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.496.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0984"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.922.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.737.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.492.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.600.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0802"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "268.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.885.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.25%  "
